# 06/11/2025 Fixed Showing up of Draft Button for IT PIC (Accepted Tix)
#
# Adds six new tracking columns (P:U) to the AutoClosedTicket report:
#   ASSIGNED IT PIC, ASSIGNED DATE TIME, RESOLVED DATE TIME,
#   SLA HOURS, ACTUAL HOURS, HIT OR MISS
# with the same header fill used by the existing header row, a
# mm/dd/yyyy hh:mm AM/PM number format on the two date/time columns,
# matching column widths, and moves the active selection to R7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells P4:U4 -------------------------------------------------
$ws.Range("P4").Value = "ASSIGNED IT PIC"
$ws.Range("Q4").Value = "ASSIGNED DATE TIME"
$ws.Range("R4").Value = "RESOLVED DATE TIME"
$ws.Range("S4").Value = "SLA HOURS"
$ws.Range("T4").Value = "ACTUAL HOURS"
$ws.Range("U4").Value = "HIT OR MISS"

# Reuse the existing header style (theme fill) from O4 for all the new headers
$ws.Range("O4").Copy()
$ws.Range("P4:U4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Q4/R4 (ASSIGNED DATE TIME / RESOLVED DATE TIME) additionally get a
# date-time display format, layered on top of the header fill above
$ws.Range("Q4:R4").NumberFormat = "mm/dd/yyyy\ hh:mm\ AM/PM"

# --- Column widths for the new columns --------------------------------------
# (ColumnWidth is expressed in "characters"; the stored xlsx width is
#  ColumnWidth + 5/6, so subtract that offset to hit the target widths)
$ws.Columns.Item(16).ColumnWidth = 26.7109375 - 0.8333333333333334   # P
$ws.Columns.Item(17).ColumnWidth = 31.7109375 - 0.8333333333333334   # Q
$ws.Columns.Item(18).ColumnWidth = 28.42578125 - 0.8333333333333334  # R
$ws.Columns.Item(19).ColumnWidth = 27.85546875 - 0.8333333333333334  # S
$ws.Columns.Item(20).ColumnWidth = 26.85546875 - 0.8333333333333334  # T
$ws.Columns.Item(21).ColumnWidth = 23.28515625 - 0.8333333333333334  # U

# --- Update the current selection/view --------------------------------------
$ws.Range("R7").Select()
